$d = $word.ActiveDocument

# --- Update date in the 'first page' header (header3.xml / rId12) ---
$dateSection = $d.Sections.First
$dateHeader = $dateSection.Headers.Item(2)
$dateHeader.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null

# --- Phase 1: insert all new 'Knärot' section paragraphs (plain text only) ---
$anchor = $d.Paragraphs.Last
$newParas = @()

# paragraph 0
$anchor.Range.InsertParagraphAfter()
$p0 = $d.Paragraphs.Last
$p0.Style = "Heading1"
$cur = $p0.Range
$cur.Collapse(0)
$cur.InsertAfter('Knärot – ekologi samt krav på livsmiljön')
$newParas += $p0
$anchor = $p0

# paragraph 1
$anchor.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Style = "Normal"
$cur = $p1.Range
$cur.Collapse(0)
$cur.InsertAfter('Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).')
$newParas += $p1
$anchor = $p1

# paragraph 2
$anchor.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Style = "Normal"
$cur = $p2.Range
$cur.Collapse(0)
$cur.InsertAfter('Samuel Johnsons doktorsavhandling ')
$cur.InsertAfter('“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“')
$cur.InsertAfter(' (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ')
$cur.InsertAfter('“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ')
$cur.InsertAfter('Vidare ')
$cur.InsertAfter('“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”')
$newParas += $p2
$anchor = $p2

# paragraph 3
$anchor.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Style = "Normal"
$cur = $p3.Range
$cur.Collapse(0)
$cur.InsertAfter('Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: ')
$cur.InsertAfter('“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”')
$newParas += $p3
$anchor = $p3

# paragraph 4
$anchor.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$p4.Style = "Normal"
$cur = $p4.Range
$cur.Collapse(0)
$cur.InsertAfter('En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).')
$newParas += $p4
$anchor = $p4

# paragraph 5
$anchor.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Last
$p5.Style = "Normal"
$cur = $p5.Range
$cur.Collapse(0)
$cur.InsertAfter('Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).')
$newParas += $p5
$anchor = $p5

# paragraph 6
$anchor.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Last
$p6.Style = "Heading2"
$cur = $p6.Range
$cur.Collapse(0)
$cur.InsertAfter('Referenser - knärot')
$newParas += $p6
$anchor = $p6

# paragraph 7
$anchor.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs.Last
$p7.Style = "Normal"
$cur = $p7.Range
$cur.Collapse(0)
$cur.InsertAfter('de Graaf M & Roberts M.R., 2009. ')
$cur.InsertAfter('Short-term response of the herbaceous layer within leave patches after harvest. ')
$cur.InsertAfter('Forest Ecology and Management 257, 1014-1025')
$newParas += $p7
$anchor = $p7

# paragraph 8
$anchor.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs.Last
$p8.Style = "Normal"
$cur = $p8.Range
$cur.Collapse(0)
$cur.InsertAfter('Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. ')
$cur.InsertAfter('Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ')
$cur.InsertAfter('Ecological Applications, 22, 2049-2064 ')
$newParas += $p8
$anchor = $p8

# paragraph 9
$anchor.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs.Last
$p9.Style = "Normal"
$cur = $p9.Range
$cur.Collapse(0)
$cur.InsertAfter('Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. ')
$cur.InsertAfter('Interactive effects of drought and edge exposure on old-growth forest understory species. ')
$cur.InsertAfter('Landscape Ecology, 37, sid 1839-1853')
$newParas += $p9
$anchor = $p9

# paragraph 10
$anchor.Range.InsertParagraphAfter()
$p10 = $d.Paragraphs.Last
$p10.Style = "Normal"
$cur = $p10.Range
$cur.Collapse(0)
$cur.InsertAfter('Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. ')
$cur.InsertAfter('Biological legacies buffer local species extinction after logging. ')
$cur.InsertAfter('Journal of Applied Ecology. 51, 53-62.')
$newParas += $p10
$anchor = $p10

# paragraph 11
$anchor.Range.InsertParagraphAfter()
$p11 = $d.Paragraphs.Last
$p11.Style = "Normal"
$cur = $p11.Range
$cur.Collapse(0)
$cur.InsertAfter('Skogsstyrelsen, 2022. ')
$cur.InsertAfter('Vägledning för hänsyn till knärot. ')
$cur.InsertAfter('https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/')
$newParas += $p11
$anchor = $p11

# paragraph 12
$anchor.Range.InsertParagraphAfter()
$p12 = $d.Paragraphs.Last
$p12.Style = "Normal"
$cur = $p12.Range
$cur.Collapse(0)
$cur.InsertAfter('SLU Artdatabanken, 2021. ')
$cur.InsertAfter('Artfaktablad. Naturvård – artfakta. ')
$cur.InsertAfter('SLU Artdatabanken, Uppsala ')
$newParas += $p12
$anchor = $p12

# --- Phase 2: apply italic formatting to designated sub-ranges ---

# italics for paragraph 2
$searchFrom = $p2.Range.Start
$f = $d.Range($searchFrom, $p2.Range.End)
$f.Find.Execute('Samuel Johnsons doktorsavhandling ') | Out-Null
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p2.Range.End)
$f.Find.Execute('“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“') | Out-Null
$f.Font.Italic = 1
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p2.Range.End)
$f.Find.Execute(' (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ') | Out-Null
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p2.Range.End)
$f.Find.Execute('“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ') | Out-Null
$f.Font.Italic = 1
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p2.Range.End)
$f.Find.Execute('Vidare ') | Out-Null
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p2.Range.End)
$f.Find.Execute('“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”') | Out-Null
$f.Font.Italic = 1
$searchFrom = $f.End

# italics for paragraph 3
$searchFrom = $p3.Range.Start
$f = $d.Range($searchFrom, $p3.Range.End)
$f.Find.Execute('Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: ') | Out-Null
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p3.Range.End)
$f.Find.Execute('“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”') | Out-Null
$f.Font.Italic = 1
$searchFrom = $f.End

# italics for paragraph 7
$searchFrom = $p7.Range.Start
$f = $d.Range($searchFrom, $p7.Range.End)
$f.Find.Execute('de Graaf M & Roberts M.R., 2009. ') | Out-Null
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p7.Range.End)
$f.Find.Execute('Short-term response of the herbaceous layer within leave patches after harvest. ') | Out-Null
$f.Font.Italic = 1
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p7.Range.End)
$f.Find.Execute('Forest Ecology and Management 257, 1014-1025') | Out-Null
$searchFrom = $f.End

# italics for paragraph 8
$searchFrom = $p8.Range.Start
$f = $d.Range($searchFrom, $p8.Range.End)
$f.Find.Execute('Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. ') | Out-Null
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p8.Range.End)
$f.Find.Execute('Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ') | Out-Null
$f.Font.Italic = 1
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p8.Range.End)
$f.Find.Execute('Ecological Applications, 22, 2049-2064 ') | Out-Null
$searchFrom = $f.End

# italics for paragraph 9
$searchFrom = $p9.Range.Start
$f = $d.Range($searchFrom, $p9.Range.End)
$f.Find.Execute('Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. ') | Out-Null
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p9.Range.End)
$f.Find.Execute('Interactive effects of drought and edge exposure on old-growth forest understory species. ') | Out-Null
$f.Font.Italic = 1
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p9.Range.End)
$f.Find.Execute('Landscape Ecology, 37, sid 1839-1853') | Out-Null
$searchFrom = $f.End

# italics for paragraph 10
$searchFrom = $p10.Range.Start
$f = $d.Range($searchFrom, $p10.Range.End)
$f.Find.Execute('Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. ') | Out-Null
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p10.Range.End)
$f.Find.Execute('Biological legacies buffer local species extinction after logging. ') | Out-Null
$f.Font.Italic = 1
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p10.Range.End)
$f.Find.Execute('Journal of Applied Ecology. 51, 53-62.') | Out-Null
$searchFrom = $f.End

# italics for paragraph 11
$searchFrom = $p11.Range.Start
$f = $d.Range($searchFrom, $p11.Range.End)
$f.Find.Execute('Skogsstyrelsen, 2022. ') | Out-Null
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p11.Range.End)
$f.Find.Execute('Vägledning för hänsyn till knärot. ') | Out-Null
$f.Font.Italic = 1
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p11.Range.End)
$f.Find.Execute('https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/') | Out-Null
$searchFrom = $f.End

# italics for paragraph 12
$searchFrom = $p12.Range.Start
$f = $d.Range($searchFrom, $p12.Range.End)
$f.Find.Execute('SLU Artdatabanken, 2021. ') | Out-Null
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p12.Range.End)
$f.Find.Execute('Artfaktablad. Naturvård – artfakta. ') | Out-Null
$f.Font.Italic = 1
$searchFrom = $f.End
$f = $d.Range($searchFrom, $p12.Range.End)
$f.Find.Execute('SLU Artdatabanken, Uppsala ') | Out-Null
$searchFrom = $f.End

Write-Output "Done. Paragraphs now: $($d.Paragraphs.Count)"
